# Replace the working set of sequence rows (columns B:E, rows 2-33) with the
# new batch of images/words/categories described in the commit "Add working
# set of sequences". Column A (the 0-based row index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B (count), C (image), D (word), E (category)
$rows = @(
    @(2,  95,  "dog/dog025.jpg",      "regnen",   "dog"),
    @(3,  125, "flower/flower014.jpg","antun",    "flower"),
    @(4,  41,  "dog/dog013.jpg",      "dauern",   "dog"),
    @(5,  23,  "flower/flower022.jpg","hupen",    "flower"),
    @(6,  10,  "flower/flower023.jpg","segeln",   "flower"),
    @(7,  82,  "flower/flower024.jpg","biegen",   "flower"),
    @(8,  98,  "dog/dog012.jpg",      "sondern",  "dog"),
    @(9,  51,  "dog/dog015.jpg",      "stechen",  "dog"),
    @(10, 50,  "flower/flower002.jpg","stärken",  "flower"),
    @(11, 27,  "dog/dog017.jpg",      "fliehen",  "dog"),
    @(12, 121, "flower/flower027.jpg","lehnen",   "flower"),
    @(13, 103, "flower/flower018.jpg","kehren",   "flower"),
    @(14, 109, "dog/dog021.jpg",      "hoffen",   "dog"),
    @(15, 54,  "flower/flower010.jpg","tauschen", "flower"),
    @(16, 40,  "flower/flower015.jpg","gelten",   "flower"),
    @(17, 45,  "dog/dog022.jpg",      "füllen",   "dog"),
    @(18, 126, "flower/flower026.jpg","pflegen",  "flower"),
    @(19, 52,  "dog/dog005.jpg",      "fühlen",   "dog"),
    @(20, 96,  "flower/flower005.jpg","strahlen", "flower"),
    @(21, 89,  "flower/flower012.jpg","schicken", "flower"),
    @(22, 70,  "dog/dog004.jpg",      "saufen",   "dog"),
    @(23, 31,  "flower/flower020.jpg","tagen",    "flower"),
    @(24, 15,  "dog/dog020.jpg",      "langen",   "dog"),
    @(25, 68,  "flower/flower028.jpg","fliegen",  "flower"),
    @(26, 114, "dog/dog002.jpg",      "gründen",  "dog"),
    @(27, 7,   "flower/flower000.jpg","starten",  "flower"),
    @(28, 115, "dog/dog029.jpg",      "fesseln",  "dog"),
    @(29, 73,  "dog/dog026.jpg",      "rücken",   "dog"),
    @(30, 94,  "dog/dog016.jpg",      "enden",    "dog"),
    @(31, 65,  "dog/dog003.jpg",      "hauen",    "dog"),
    @(32, 42,  "dog/dog027.jpg",      "rasen",    "dog"),
    @(33, 71,  "flower/flower009.jpg","krachen",  "flower")
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    $ws.Range("D$rowNum").Value = $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
}
